# Update the 10.b.1 indicator title (both language variants) to 10.b.1.1
# in the header row (A1 = Kyrgyz, C1 = English). B1 already holds the
# Russian title which already included the ".1.1" numbering and is
# unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "10.b.1.1 Агымдардын түрлөрү жана алуучу өлкөлөр жана донор-өлкөлөр боюнча бөлунүшүндөгү  өнүктүрүү максатында ресурстар агымынын жалпы көлөмү"
$ws.Range("C1").Value = "10.b.1.1 Total resource flows for development, by recipient and donor countries and type of flow (e.g. official development assistance, foreign direct investment and other flows)"

# Leave the selection on L8 (matches the cursor position recorded when the
# workbook was last saved).
$ws.Range("L8").Select()
